$d = $word.ActiveDocument
$s = $d.Shapes.Item(2)
$tf = $s.TextFrame
$cr = $tf.ContainingRange
Write-Host "ContainingRange text: [$($cr.Text)]"
Write-Host "ContainingRange Start/End: $($cr.Start) / $($cr.End)"
